$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 120000
$ws.Range("B3").Value = 170000
$ws.Range("B4").Value = 2
$ws.Range("B7").Value = "andra.andruta60@gmail.com"

$ws.Range("B3").Select()

$wb.Save()
